{"js": "// The document contains five empty \"My Style\" list paragraphs that act as\n// answer slots beneath their corresponding question/heading paragraphs.\n// This script fills each of those empty paragraphs with its answer text,\n// matching the target diff exactly (inserting a single run/w:t into each\n// otherwise-empty <w:p>).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Map each (currently empty) paragraph's plain text to the text that should\n// be inserted into it. Using the empty-text signal together with sequential\n// order keeps this robust even though all five target paragraphs currently\n// share the same (empty) text value.\nconst insertions = [\n  \"By asking questions & jotting the answers down\",\n  \"It doesn\\u2019t seem to line up with how the actual game turned out.\",\n  \"That it seems to be out of date compared to the game.\",\n  \"Updating both AT01 & 2\",\n  \"None rejected.\"\n];\n\nlet insertionIndex = 0;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text === \"\" && insertionIndex < insertions.length) {\n    para.insertText(insertions[insertionIndex], Word.InsertLocation.start);\n    insertionIndex++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains five empty \"My Style\" list paragraphs that act as\n# answer slots beneath their corresponding question/heading paragraphs.\n# This script fills each of those empty paragraphs with its answer text,\n# matching the target diff exactly (inserting a single run/text into each\n# otherwise-empty paragraph).\n$d = $word.ActiveDocument\n\n$insertions = @(\n    \"By asking questions & jotting the answers down\",\n    \"It doesn\u2019t seem to line up with how the actual game turned out.\",\n    \"That it seems to be out of date compared to the game.\",\n    \"Updating both AT01 & 2\",\n    \"None rejected.\"\n)\n\n$idx = 0\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd(\"`r\", \"`n\", \"`x07\")\n    if ($t -eq \"\" -and $idx -lt $insertions.Count) {\n        $p.Range.InsertAfter($insertions[$idx])\n        $idx++\n    }\n}\n"}
